$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.194.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6186"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07358"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2916"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.838.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.973"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6683"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008937"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.851"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.176.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.076.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "235.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.350"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9988"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1396"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.544"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.492"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05778"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.107"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.086"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.210"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.851"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7288"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.141"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.607"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.222.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01753"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.248"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9043"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.981.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000118"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.125"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4022"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1132"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.98%  "

Write-Host "Applied all cell updates"